# RPA datasets push 2024-05-02
# Insert two new IPO dataset rows at the top of the data table (rows 2-3),
# pushing all existing rows down by two.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 blank rows before the current row 2 (shifts rows 2.. down to 4..)
$ws.Rows("2:3").Insert()

# Row insertion copies the formatting of the row above (the bold, centered
# header row). Strip that back to the plain, unstyled look the rest of the
# data rows use.
$ws.Range("A2:Y3").ClearFormats()

# Force the text-bearing columns to plain "Text" format so Excel doesn't
# auto-convert date-looking / percent-looking strings into real
# dates/numbers (matches the source data, where these are shared strings).
$ws.Range("A2:E3").NumberFormat = "@"
$ws.Range("N2:O3").NumberFormat = "@"
$ws.Range("Y2:Y3").NumberFormat = "@"

# New row 2: 디앤디파마텍 (D&D Pharmatech)
$ws.Range("A2").Value = "2024-04-12"
$ws.Range("B2").Value = "2024-04-18"
$ws.Range("C2").Value = "2024-05-02"
$ws.Range("D2").Value = "한국"
$ws.Range("E2").Value = "디앤디파마텍"
$ws.Range("F2").Value = 1100000
$ws.Range("G2").Value = 1100000
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 22000
$ws.Range("J2").Value = 26000
$ws.Range("K2").Value = 10429232
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 33000
$ws.Range("N2").Value = "848.50:1"
$ws.Range("O2").Value = "10.96%"
$ws.Range("P2").Value = -75676750274
$ws.Range("Q2").Value = -68652978862
$ws.Range("R2").Value = -9506668082
$ws.Range("S2").Value = -69862474811
$ws.Range("T2").Value = -137025491259
$ws.Range("U2").Value = 3014576074
$ws.Range("V2").Value = 0
$ws.Range("W2").Value = 0
$ws.Range("X2").Value = 0
$ws.Range("Y2").Value = "대사성질환 치료제 등"

# New row 3: 유안타제16호스팩 (Yuanta No.16 SPAC)
$ws.Range("A3").Value = "2024-04-15"
$ws.Range("B3").Value = "2024-04-16"
$ws.Range("C3").Value = "2024-05-02"
$ws.Range("D3").Value = "유안타"
$ws.Range("E3").Value = "유안타제16호스팩"
$ws.Range("F3").Value = 5150000
$ws.Range("G3").Value = 5150000
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 2000
$ws.Range("J3").Value = 2000
$ws.Range("K3").Value = 5510000
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = 2000
$ws.Range("N3").Value = "1,050.42:1"
$ws.Range("O3").Value = "-"
$ws.Range("P3").Value = 0
$ws.Range("Q3").Value = 0
$ws.Range("R3").Value = 0
$ws.Range("S3").Value = 0
$ws.Range("T3").Value = 0
$ws.Range("U3").Value = 0
$ws.Range("V3").Value = 0
$ws.Range("W3").Value = 0
$ws.Range("X3").Value = 0
$ws.Range("Y3").Value = "금융 지원 서비스(기업인수목적회사)"
